$wb = $excel.ActiveWorkbook

# --- "ir" sheet: move the SKIP flag from row 2 (99CENTS) to row 3 (USFOOD)
#     and update the remembered selection to E4.
$wsIr = $wb.Worksheets.Item("ir")
$wsIr.Range("E2").ClearContents()
$wsIr.Range("E3").Value = $true
$wsIr.Range("E4").Select() | Out-Null

# --- "simpleton" sheet: check the SKIP flag on row 2 (USFOOD)
#     and update the remembered selection to H3. This sheet was the
#     previously active tab, so it loses tabSelected once another
#     sheet is activated below.
$wsSimpleton = $wb.Worksheets.Item("simpleton")
$wsSimpleton.Range("H2").Value = $true
$wsSimpleton.Range("H3").Select() | Out-Null

# --- "tier" sheet becomes the active tab (activeTab goes from 3 to 1).
$wsTier = $wb.Worksheets.Item("tier")
$wsTier.Range("E16").Select() | Out-Null
